# Scheduled-runner style refresh of raw market-board figures (columns H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets. Values are plain
# numeric literals (no formulas), so each changed cell is just re-written
# with its refreshed figure; a few rows also gain/lose a cell (M/N) where
# the source feed now does/doesn't report that figure.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3368.2307
$ws.Range("I40").Value = 1933.6666
$ws.Range("K40").Value = 1933.6666
$ws.Range("M40").Value = -1758.6666

$ws.Range("H74").Value = 8229.333000000001
$ws.Range("I74").Value = 7844
$ws.Range("K74").Value = 7844
$ws.Range("M74").Value = -6908

$ws.Range("H77").Value = 8229.333000000001
$ws.Range("I77").Value = 7844
$ws.Range("K77").Value = 39220
$ws.Range("M77").Value = -34540

$ws.Range("H100").Value = 13660166
$ws.Range("I100").Value = 13894021
$ws.Range("K100").Value = 13894021
$ws.Range("M100").Value = -13893480

$ws.Range("H132").Value = 4369.3706
$ws.Range("I132").Value = 4709.857
$ws.Range("J132").Value = 3177.6667
$ws.Range("K132").Value = 14129.571
$ws.Range("L132").Value = 9533.000100000001
$ws.Range("M132").Value = -11599.571
$ws.Range("N132").Value = -14593.0001

$ws.Range("H137").Value = 5209862
$ws.Range("I137").Value = 1146440.6
$ws.Range("J137").Value = 7361085
$ws.Range("K137").Value = 3439321.8
$ws.Range("L137").Value = 22083255
$ws.Range("M137").Value = -3436771.8
$ws.Range("N137").Value = -22088355

$ws.Range("H141").Value = 5344.1
$ws.Range("I141").Value = 4180.75
$ws.Range("K141").Value = 12542.25
$ws.Range("M141").Value = -7362.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 224382.2
$ws.Range("I2").Value = 55705.5
$ws.Range("K2").Value = 55705.5
$ws.Range("M2").Value = -55592.5

$ws.Range("H32").Value = 5347.4
$ws.Range("I32").Value = 5531.591
$ws.Range("K32").Value = 5531.591
$ws.Range("M32").Value = -5244.591

$ws.Range("H116").Value = 224382.2
$ws.Range("I116").Value = 55705.5
$ws.Range("K116").Value = 55705.5
$ws.Range("M116").Value = -53411.5

$ws.Range("H132").Value = 3295.2307
$ws.Range("I132").Value = 3008.1365
$ws.Range("K132").Value = 9024.4095
$ws.Range("M132").Value = -6494.4095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 224382.2
$ws.Range("I3").Value = 55705.5
$ws.Range("K3").Value = 55705.5
$ws.Range("M3").Value = -55591.5

$ws.Range("H20").Value = 4190.385
$ws.Range("I20").Value = 1830.2
$ws.Range("K20").Value = 1830.2
$ws.Range("M20").Value = -1583.2

$ws.Range("H41").Value = 193217.5
$ws.Range("J41").Value = 193217.5
$ws.Range("L41").Value = 193217.5
$ws.Range("N41").Value = -193993.5

$ws.Range("H48").Value = 198905
$ws.Range("J48").Value = 198905
$ws.Range("L48").Value = 198905
$ws.Range("N48").Value = -199735

$ws.Range("H99").Value = 12171.625
$ws.Range("I99").Value = 13238.839
$ws.Range("J99").Value = 8495.666999999999
$ws.Range("K99").Value = 13238.839
$ws.Range("L99").Value = 8495.666999999999
$ws.Range("M99").Value = -11740.839
$ws.Range("N99").Value = -11491.667

$ws.Range("H105").Value = 94664.664
$ws.Range("I105").Value = 123886.336
$ws.Range("K105").Value = 123886.336
$ws.Range("M105").Value = -122139.336

$ws.Range("H134").Value = 13822.385
$ws.Range("I134").Value = 15153.728
$ws.Range("K134").Value = 45461.18399999999
$ws.Range("M134").Value = -42926.18399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2300
$ws.Range("I31").Value = 749.38464
$ws.Range("K31").Value = 749.38464
$ws.Range("M31").Value = -454.38464

$ws.Range("H33").Value = 6795.6
$ws.Range("I33").Value = 3494.5
$ws.Range("K33").Value = 3494.5
$ws.Range("M33").Value = -3115.5

$ws.Range("H34").Value = 2300
$ws.Range("I34").Value = 749.38464
$ws.Range("K34").Value = 749.38464
$ws.Range("M34").Value = -547.38464

$ws.Range("H132").Value = 41733240
$ws.Range("I132").Value = 83364104
$ws.Range("K132").Value = 250092312
$ws.Range("M132").Value = -250089782

$ws.Range("H134").Value = 2321432.2
$ws.Range("I134").Value = 2983580.2
$ws.Range("J134").Value = 3913.6667
$ws.Range("K134").Value = 8950740.600000001
$ws.Range("L134").Value = 11741.0001
$ws.Range("M134").Value = -8948205.600000001
$ws.Range("N134").Value = -16811.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 218.57895
$ws.Range("I2").Value = 159.4
$ws.Range("J2").Value = 284.33334
$ws.Range("K2").Value = 956.4000000000001
$ws.Range("L2").Value = 1706.00004
$ws.Range("M2").Value = -843.4000000000001
$ws.Range("N2").Value = -1932.00004

$ws.Range("H38").Value = 1178.2903
$ws.Range("I38").Value = 213.5
$ws.Range("J38").Value = 1637.7142
$ws.Range("K38").Value = 640.5
$ws.Range("L38").Value = 4913.142599999999
$ws.Range("M38").Value = -293.5
$ws.Range("N38").Value = -5607.142599999999

$ws.Range("H68").Value = 18521208
$ws.Range("J68").Value = 25002990
$ws.Range("L68").Value = 75008970
$ws.Range("N68").Value = -75010592

$ws.Range("H71").Value = 18521208
$ws.Range("J71").Value = 25002990
$ws.Range("L71").Value = 225026910
$ws.Range("N71").Value = -225035022

$ws.Range("H134").Value = 2580.5715
$ws.Range("I134").Value = 2163.6924
$ws.Range("K134").Value = 6491.0772
$ws.Range("M134").Value = -1421.0772

$ws.Range("H139").Value = 3003451.5
$ws.Range("I139").Value = 6001904
$ws.Range("K139").Value = 18005712
$ws.Range("M139").Value = -18000572

$ws.Range("H140").Value = 8631.454
$ws.Range("I140").Value = 9069.700000000001
$ws.Range("K140").Value = 27209.1
$ws.Range("M140").Value = -22029.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6400.091
$ws.Range("I70").Value = 6758.706
$ws.Range("J70").Value = 5180.8
$ws.Range("K70").Value = 6758.706
$ws.Range("L70").Value = 5180.8
$ws.Range("M70").Value = -6488.706
$ws.Range("N70").Value = -5720.8

$ws.Range("H73").Value = 6400.091
$ws.Range("I73").Value = 6758.706
$ws.Range("J73").Value = 5180.8
$ws.Range("K73").Value = 6758.706
$ws.Range("L73").Value = 5180.8
$ws.Range("M73").Value = -5822.706
$ws.Range("N73").Value = -7052.8

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H122").Value = 11565.875
$ws.Range("I122").Value = 12182.429
$ws.Range("K122").Value = 36547.287
$ws.Range("M122").Value = -34097.287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1306
$ws.Range("J16").Value = 825
$ws.Range("L16").Value = 825
$ws.Range("N16").Value = -1165

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 4589.357
$ws.Range("I122").Value = 4029.25
$ws.Range("K122").Value = 12087.75
$ws.Range("M122").Value = -9637.75

$ws.Range("H132").Value = 13535.471
$ws.Range("I132").Value = 21660.777
$ws.Range("K132").Value = 64982.33099999999
$ws.Range("M132").Value = -62452.33099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 40913.23
$ws.Range("I100").Value = 21024.727
$ws.Range("K100").Value = 42049.454
$ws.Range("M100").Value = -41508.454
